# Saturday shift time changed from 09:00-13:00 to 10:00-14:00.
# This touches four sheets:
#   - Pianificazione: the "Sabato - Fascia" column (H) for the affected rows
#   - Assegnazioni: the per-assignment Sabato rows (E/F/G/H/I/M columns)
#   - Copertura: the Sabato coverage numbers that move from the old slot
#     (09:00-10:00) to the new one (13:00-14:00)
#   - Warnings: the warning that is no longer applicable (demand fully
#     covered at 13:00-13:45 now) is removed, shifting the remaining rows up

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Pianificazione sheet - update the Sabato slot for the rows that had
#    the old 09:00-13:00 shift.
# ---------------------------------------------------------------------
$wsPian = $wb.Worksheets.Item("Pianificazione")
$pianRows = @(2, 5, 8, 9, 19, 21)
foreach ($r in $pianRows) {
    $wsPian.Range("H$r").Value = "10:00-14:00"
}

# ---------------------------------------------------------------------
# 2) Assegnazioni sheet - update the detailed assignment columns for the
#    rows belonging to the Sabato 09:00-13:00 shift.
# ---------------------------------------------------------------------
$wsAsseg = $wb.Worksheets.Item("Assegnazioni")
$assegRows = @(96, 97, 98, 99, 100, 101)
foreach ($r in $assegRows) {
    $wsAsseg.Range("E$r").Value = "10:00-14:00"
    $wsAsseg.Range("F$r").Value = "10:00"
    $wsAsseg.Range("G$r").Value = "14:00"
    $wsAsseg.Range("H$r").Value = "10:00"
    $wsAsseg.Range("I$r").Value = "14:00"
    $wsAsseg.Range("M$r").Value = "AUTO_10:00-14:00_240"
}

# ---------------------------------------------------------------------
# 3) Copertura sheet - the coverage figures shift together with the
#    employees: the old 09:00-10:00 slots lose their coverage (Coperta ->
#    0, Gap -> 0) while the 13:00-14:00 slots gain it (Coperta -> 6,
#    Gap -> 1).
# ---------------------------------------------------------------------
$wsCop = $wb.Worksheets.Item("Copertura")

$copRowsVacated = @(512, 513, 514, 515)
foreach ($r in $copRowsVacated) {
    $wsCop.Range("D$r").Value = 0
    $wsCop.Range("E$r").Value = 0
}

$copRowsFilled = @(528, 529, 530, 531)
foreach ($r in $copRowsFilled) {
    $wsCop.Range("D$r").Value = 6
    $wsCop.Range("E$r").Value = 1
}

# ---------------------------------------------------------------------
# 4) Warnings sheet - the "Domanda 7.0 ... Slot Sab 13:00 ..." warning no
#    longer applies, so its row is removed (remaining rows shift up).
# ---------------------------------------------------------------------
$wsWarn = $wb.Worksheets.Item("Warnings")
$wsWarn.Rows.Item(5).Delete()
